# edit.ps1 - applies the draft-gandhi-ippm-stamp-srpm-01.pptx revision:
#   - "message"/"LM message" -> "test packet"/"direct measurement packet"
#     wording tweaks on slides 4, 9, 11, 16 (STAMP -> direct measurement
#     terminology cleanup)
#   - shrink the height of the "Destination Node Address TLV" callout box
#     on slide 9 to match its (now shorter) text

$p = $ppt.ActivePresentation

# NOTE: this interpreter's function-parameter binding only works reliably
# with *positional* args -- passing COM objects / strings via `-Name value`
# silently binds nothing, so every call below uses positional arguments.
function Replace-InShapeText {
    param($Shape, [string]$OldText, [string]$NewText)
    $tr = $Shape.TextFrame.TextRange
    $full = $tr.Text
    $idx = $full.IndexOf($OldText)
    if ($idx -lt 0) {
        throw "Text not found in shape '$($Shape.Name)': $OldText"
    }
    $sub = $tr.Characters($idx + 1, $OldText.Length)
    $sub.Text = $NewText
}

# ---- Slide 4 : "Review Comments" bullet list ----------------------------
$s4 = $p.Slides.Item(4)
$sh4 = $s4.Shapes.Item(2)   # "Content Placeholder 2"

Replace-InShapeText $sh4 `
    "Updates RFC 8762 due to new field (control code) in the message" `
    "Updates RFC 8762 due to new field (control code) in the test packet"

Replace-InShapeText $sh4 `
    "Indicate new packet loss message is for direct measurement" `
    "Indicate new packet loss packet is for direct measurement"

Replace-InShapeText $sh4 `
    "Move Receive Counter and other Reply message fields to Section 4.1 from 3.2" `
    "Move Receive Counter and other Reply test packet fields to Section 4.1 from 3.2"

# ---- Slide 9 : "Destination Node Address TLV" callout --------------------
$s9 = $p.Slides.Item(9)
$sh9 = $s9.Shapes.Item(5)   # "Rectangle 8"

Replace-InShapeText $sh9 `
    "Indicates the address of the intended recipient node of the test packet message.  " `
    "Indicates the address of the intended recipient node of the test packet.  "

# Box shrinks vertically now that the text is one line shorter.
$sh9.Height = 2772234 / 12700

# ---- Slide 11 : "Hardware Implementation considerations" -----------------
$s11 = $p.Slides.Item(11)
$sh11 = $s11.Shapes.Item(4)   # "Content Placeholder 6"

Replace-InShapeText $sh11 `
    "Separate UDP port + LM message format eliminate the complexity in Hardware" `
    "Separate UDP port + direct measurement packet format eliminate the complexity in Hardware"

Replace-InShapeText $sh11 `
    "Counter at fixed location (Eth 18, IPv6 40, UDP 8, Seq 4, Total = 70 Byte)" `
    "Counter at fixed location (Eth 18, IPv6 40, UDP 8, Seq 4, Total = 70 Byte), not deeper in the packet"

# ---- Slide 16 : Session-Sender / Session-Reflector packet diagrams -------
$s16 = $p.Slides.Item(16)

$shSender = $s16.Shapes.Item(2)     # "Rectangle 4"
Replace-InShapeText $shSender `
    "`t    Figure: Session-Sender Message Format" `
    "`t Figure: Session-Sender Test Packet Format"

$shReflector = $s16.Shapes.Item(3)  # "Rectangle 13"
Replace-InShapeText $shReflector `
    "`t   Figure: Session-Reflector Message Format" `
    "`tFigure: Session-Reflector Test Packet Format"

Write-Output "Done."
